# Merge the split "<id>...</id>" runs (tag / text / tag) into a single run
# with the text "<id>pNNNv_N</id>", for each of the six id markers in the
# document. Word's Find & Replace naturally collapses the matched range
# into one run using the formatting of the start of the match, which is
# exactly the tag-colored (brown / Courier New) formatting we want to keep.

$d = $word.ActiveDocument

$ids = @("p020v_1", "p020v_2", "p020v_3", "p020v_4", "p020v_5", "p020v_6")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $rng = $d.Content
    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2)
}

Write-Output "done"
